$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Reference cells already carrying the pre-existing cell styles we want to
# replicate onto the new rows (copy/paste-format keeps the shared style table
# tidy instead of minting near-duplicate font/xf entries one property at a time):
#   style index 3 (wrap-text Arial body cell, column B)  -> B70
#   style index 4 (wrap-text hyperlink cell, column B)   -> B11
# Column A keeps the plain default style ("1") that the column already carries,
# so those cells need no extra formatting at all.
$styleRefBody = $ws.Range("B70")
$styleRefLink = $ws.Range("B11")

$rows = @(
    ,@(72, 16.8, 'Oliver, Pedro Juan', 'http://viaf.org/viaf/34607051', 0)
    ,@(73, 16.8, 'Gronovius, Abraham', 'http://viaf.org/viaf/61677396', 0)
    ,@(74, 16.8, 'Bleyswick, Francois van', 'http://viaf.org/viaf/12384559', 0)
    ,@(75, 16.8, 'Barbier, Jules', 'http://viaf.org/viaf/86596074', 0)
    ,@(76, 16.8, 'Bardili, Guilielmus Henricus', 'http://viaf.org/viaf/49966518', 0)
    ,@(77, 16.8, 'Staveren, Augustinus van', 'http://viaf.org/viaf/89208065', 0)
    ,@(78, 16.8, 'Baret, Eugene', 'http://viaf.org/viaf/36964332', 0)
    ,@(79, 16.8, 'Barker, Edmund Henry', 'http://viaf.org/viaf/39716436', 0)
    ,@(80, 16.8, 'Baron, Augustus', 'http://viaf.org/viaf/7193500', 0)
    ,@(81, 33.6, 'Barreau, H.', 'http://viaf.org/viaf/7691155044886272520003', 0)
    ,@(82, 16.8, 'Peronne, Joseph-Max', 'http://viaf.org/viaf/61854282', 0)
    ,@(83, 16.8, 'Ecalle, Pierre Felix', 'http://viaf.org/viaf/56646541', 0)
    ,@(84, 16.8, 'Charpentier, Jean-Pierre', 'http://viaf.org/viaf/22274769', 0)
    ,@(85, 16.8, 'Vincent, C.J.B.J.', 'http://viaf.org/viaf/121850838', 0)
    ,@(86, 16.8, 'Barth, Friedrich Gottlieb', 'http://viaf.org/viaf/908603', 0)
    ,@(87, 16.8, 'Vinet, Elie', 'http://viaf.org/viaf/27100990', 0)
    ,@(88, 16.8, 'Accorso, Mariangelo', 'http://viaf.org/viaf/42899404', 0)
    ,@(89, 16.8, 'Rivinus, Andreas', 'http://viaf.org/viaf/40246385', 0)
    ,@(90, 16.8, 'Schrijver, Pieter', 'http://viaf.org/viaf/24732263', 0)
    ,@(91, 16.8, 'Baumstark, Anton', 'http://viaf.org/viaf/59164719', 0)
    ,@(92, 16.8, 'Baune, Jaques de la', 'http://viaf.org/viaf/39454999', 0)
    ,@(93, 16.8, 'Baxter, William', 'http://viaf.org/viaf/5266955', 0)
    ,@(94, 16.8, 'Becher, Christian', 'http://viaf.org/viaf/42743606', 0)
    ,@(95, 16.8, 'Becher, Friedrich Liebgott', 'http://viaf.org/viaf/49657969', 0)
    ,@(96, 16.8, 'Beck, Charles', 'http://viaf.org/viaf/37658375', 0)
    ,@(97, 16.8, 'Beck, Christian Daniel', 'http://viaf.org/viaf/57357633', 0)
    ,@(98, 16.8, 'Becker, Gustavus', 'http://viaf.org/viaf/42580693', 0)
    ,@(99, 16.8, 'Becker, Ulrich Justus Heinrich', 'http://viaf.org/viaf/25344923', 0)
    ,@(100, 16.8, 'Becker, William Adolf', 'http://viaf.org/viaf/88728074', 0)
    ,@(101, 16.8, 'Beesly, A.H.', 'http://viaf.org/viaf/29898952', 0)
    ,@(102, 16.8, 'Beger, Laurentius', 'http://viaf.org/viaf/56907766', 0)
    ,@(103, 16.8, 'Beier, Karl', 'http://viaf.org/viaf/444806', 0)
    ,@(104, 14.4, 'Benecke, Karl', 'http://viaf.org/viaf/263942005', 1)
    ,@(105, 16.8, 'Benedict, Traugott Friedrich', 'http://viaf.org/viaf/12739365', 0)
)

foreach ($row in $rows) {
    $r = $row[0]
    $ht = $row[1]
    $aText = $row[2]
    $bText = $row[3]
    $isLink = $row[4]

    $aCell = $ws.Cells.Item($r, 1)
    $bCell = $ws.Cells.Item($r, 2)

    $aCell.Value = $aText
    $bCell.Value = $bText

    if ($isLink -eq 1) {
        [void]$ws.Hyperlinks.Add($bCell, $bText)
        $styleRefLink.Copy()
    } else {
        $styleRefBody.Copy()
    }
    $bCell.PasteSpecial(-4122)
    $excel.CutCopyMode = $false

    $ws.Rows.Item($r).RowHeight = $ht
}

[void]$ws.Range("B106").Select()
